# Update "想去人数" (interest count) values by +1 for three events,
# on both the "展览" sheet and the "全部类型" sheet (which mirrors it).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1561
    $ws.Range("F6").Value = 23
    $ws.Range("F10").Value = 405
}
